# "Creado grafico de tipos de modelo"
#
# The sheet gained a new "MAE" metric column, inserted right before the
# existing "Tipo" column (which shifts from D -> E). B2/C2 (MSE/R2) were
# also recomputed, and the new D2 cell holds the MAE value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at D, pushing the old D ("Tipo"/"single") to E.
$ws.Columns.Item(4).Insert()

# Give the new header cell the same look as the other header cells
# (bold font, border, centered) by copying C1's formatting onto D1.
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial()

# New header + value for the MAE column.
$ws.Range("D1").Value = "MAE"
$ws.Range("D2").Value = 0.2400258936083209

# Updated MSE / R2 values.
$ws.Range("B2").Value = 0.1118655028770732
$ws.Range("C2").Value = 0.9979564821829493
